$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") rows 2-37 all move from serial date 45663 to 45664
$ws.Range("C2:C37").Value = 45664
